# Natmi following Dr Hou advice
# Update NATMI LR-pair computed values for Fn1-Itgb1 (ligand-expressing / receptor-expressing
# cell counts changed from 1 to 3, with corresponding recalculated expression/specificity values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 3164.756897753223
$ws.Range("R2").Value = 28482.81207977901
$ws.Range("S2").Value = 0.02302972545954165
$ws.Range("T2").Value = 0.02302972545954165
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 2746.613047072478
$ws.Range("R3").Value = 24719.5174236523
$ws.Range("S3").Value = 0.01998692046854545
$ws.Range("T3").Value = 0.01998692046854545
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 3840.248958923539
$ws.Range("R4").Value = 34562.24063031186
$ws.Range("S4").Value = 0.02794523626225017
$ws.Range("T4").Value = 0.02794523626225018
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 40446.68922645997
$ws.Range("R5").Value = 364020.2030381398
$ws.Range("S5").Value = 0.2943278674245289
$ws.Range("T5").Value = 0.294327867424529
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 35102.66599597228
$ws.Range("R6").Value = 315923.9939637505
$ws.Range("S6").Value = 0.2554397658029108
$ws.Range("T6").Value = 0.2554397658029108
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 49079.71171627369
$ws.Range("R7").Value = 441717.4054464632
$ws.Range("S7").Value = 0.3571497979075956
$ws.Range("T7").Value = 0.3571497979075956
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 986.5380338812028
$ws.Range("R8").Value = 8878.842304930826
$ws.Range("S8").Value = 0.007178971658710861
$ws.Range("T8").Value = 0.007178971658710863
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 856.1915884329104
$ws.Range("R9").Value = 7705.724295896193
$ws.Range("S9").Value = 0.006230449244419762
$ws.Range("T9").Value = 0.006230449244419763
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 1197.106690956388
$ws.Range("R10").Value = 10773.96021860749
$ws.Range("S10").Value = 0.008711265771496777
$ws.Range("T10").Value = 0.008711265771496779
